# Accepting either comma delimiter for MVplate
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H6: was blank; now a text value allowing a comma delimiter.
$ws.Range("H6").Value = ", s"

# G22:G27: were formulas (=F22/2 etc, each showing 64); now a text value
# allowing a comma delimiter, replacing the shared formula entirely.
$ws.Range("G22:G27").Value = "64,5"

# C2: was a plain number (1000); now a text value allowing a comma delimiter.
$ws.Range("C2").Value = "1000,5"

# Move the active selection from H6 to C2.
$ws.Range("C2").Select()
